# Auto-generated Excel COM-interop script
# Refreshes FFXIV market-price-derived columns (H:N) across several Sheets
# to reflect newly pulled Universalis price data (scheduled runner update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 980.1429000000001
$ws.Range("I28").Value = 893.3333
$ws.Range("J28").Value = 1501
$ws.Range("K28").Value = 893.3333
$ws.Range("L28").Value = 1501
$ws.Range("M28").Value = -408.3333
$ws.Range("N28").Value = -2471
$ws.Range("H100").Value = 40001370
$ws.Range("I100").Value = 50001430
$ws.Range("K100").Value = 50001430
$ws.Range("M100").Value = -50000889
$ws.Range("H132").Value = 25642562
$ws.Range("I132").Value = 27028106
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 81084318
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -81081788
$ws.Range("N132").Value = -35060
$ws.Range("H138").Value = 2806.15
$ws.Range("I138").Value = 1101.5
$ws.Range("J138").Value = 2995.5557
$ws.Range("K138").Value = 3304.5
$ws.Range("L138").Value = 8986.667099999999
$ws.Range("M138").Value = 1835.5
$ws.Range("N138").Value = -19266.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 417
$ws.Range("I2").Value = 334
$ws.Range("J2").Value = 624.5
$ws.Range("K2").Value = 334
$ws.Range("L2").Value = 624.5
$ws.Range("M2").Value = -221
$ws.Range("N2").Value = -850.5
$ws.Range("H32").Value = 4462.098
$ws.Range("I32").Value = 3399.7856
$ws.Range("K32").Value = 3399.7856
$ws.Range("M32").Value = -3112.7856
$ws.Range("H116").Value = 417
$ws.Range("I116").Value = 334
$ws.Range("J116").Value = 624.5
$ws.Range("K116").Value = 334
$ws.Range("L116").Value = 624.5
$ws.Range("M116").Value = 1960
$ws.Range("N116").Value = -5212.5
$ws.Range("H133").Value = 20950.4
$ws.Range("J133").Value = 20950.4
$ws.Range("L133").Value = 20950.4
$ws.Range("N133").Value = -26010.4
$ws.Range("H137").Value = 40684.668
$ws.Range("J137").Value = 40684.668
$ws.Range("L137").Value = 40684.668
$ws.Range("N137").Value = -50884.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 417
$ws.Range("I3").Value = 334
$ws.Range("J3").Value = 624.5
$ws.Range("K3").Value = 334
$ws.Range("L3").Value = 624.5
$ws.Range("M3").Value = -220
$ws.Range("N3").Value = -852.5
$ws.Range("H58").Value = 12500
$ws.Range("J58").Value = 12500
$ws.Range("L58").Value = 12500
$ws.Range("N58").Value = -13088
$ws.Range("H59").Value = 35666.332
$ws.Range("J59").Value = 35666.332
$ws.Range("L59").Value = 35666.332
$ws.Range("N59").Value = -37360.332
$ws.Range("H99").Value = 4310.1665
$ws.Range("I99").Value = 1422.2
$ws.Range("J99").Value = 6373
$ws.Range("K99").Value = 1422.2
$ws.Range("L99").Value = 6373
$ws.Range("M99").Value = 75.79999999999995
$ws.Range("N99").Value = -9369
$ws.Range("H134").Value = 2158.077
$ws.Range("I134").Value = 1630.7812
$ws.Range("J134").Value = 4568.5713
$ws.Range("K134").Value = 4892.3436
$ws.Range("L134").Value = 13705.7139
$ws.Range("M134").Value = -2357.3436
$ws.Range("N134").Value = -18775.7139
$ws.Range("H137").Value = 45709.5
$ws.Range("J137").Value = 45709.5
$ws.Range("L137").Value = 45709.5
$ws.Range("N137").Value = -55909.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2956.1936
$ws.Range("I31").Value = 1001.875
$ws.Range("J31").Value = 5040.8
$ws.Range("K31").Value = 1001.875
$ws.Range("L31").Value = 5040.8
$ws.Range("M31").Value = -706.875
$ws.Range("N31").Value = -5630.8
$ws.Range("H34").Value = 2956.1936
$ws.Range("I34").Value = 1001.875
$ws.Range("J34").Value = 5040.8
$ws.Range("K34").Value = 1001.875
$ws.Range("L34").Value = 5040.8
$ws.Range("M34").Value = -799.875
$ws.Range("N34").Value = -5444.8
$ws.Range("H58").Value = 3294.5085
$ws.Range("I58").Value = 1919.4667
$ws.Range("J58").Value = 7714.2856
$ws.Range("K58").Value = 1919.4667
$ws.Range("L58").Value = 7714.2856
$ws.Range("M58").Value = -1716.4667
$ws.Range("N58").Value = -8120.2856
$ws.Range("H99").Value = 12504107
$ws.Range("I99").Value = 18183428
$ws.Range("K99").Value = 18183428
$ws.Range("M99").Value = -18181930
$ws.Range("H122").Value = 3904.1667
$ws.Range("I122").Value = 1370.3334
$ws.Range("J122").Value = 6438
$ws.Range("K122").Value = 4111.0002
$ws.Range("L122").Value = 19314
$ws.Range("M122").Value = -1661.0002
$ws.Range("N122").Value = -24214
$ws.Range("H126").Value = 12504107
$ws.Range("I126").Value = 18183428
$ws.Range("K126").Value = 54550284
$ws.Range("M126").Value = -54547814
$ws.Range("H134").Value = 1443.4117
$ws.Range("I134").Value = 755.48
$ws.Range("J134").Value = 3354.3333
$ws.Range("K134").Value = 2266.44
$ws.Range("L134").Value = 10062.9999
$ws.Range("M134").Value = 268.5599999999999
$ws.Range("N134").Value = -15132.9999
$ws.Range("H136").Value = 3294.5085
$ws.Range("I136").Value = 1919.4667
$ws.Range("J136").Value = 7714.2856
$ws.Range("K136").Value = 5758.4001
$ws.Range("L136").Value = 23142.8568
$ws.Range("M136").Value = -3208.4001
$ws.Range("N136").Value = -28242.8568
$ws.Range("H137").Value = 45471.43
$ws.Range("J137").Value = 45471.43
$ws.Range("L137").Value = 45471.43
$ws.Range("N137").Value = -55671.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1442.8572
$ws.Range("J36").Value = 1980
$ws.Range("L36").Value = 5940
$ws.Range("N36").Value = -6278
$ws.Range("H113").Value = 3572028.8
$ws.Range("I113").Value = 618.05884
$ws.Range("J113").Value = 6945028
$ws.Range("K113").Value = 1854.17652
$ws.Range("L113").Value = 20835084
$ws.Range("M113").Value = 315.82348
$ws.Range("N113").Value = -20839424
$ws.Range("H122").Value = 3091.5
$ws.Range("I122").Value = 1154
$ws.Range("K122").Value = 10386
$ws.Range("M122").Value = -7936
$ws.Range("H131").Value = 693.86
$ws.Range("I131").Value = 261.85
$ws.Range("J131").Value = 801.8625
$ws.Range("K131").Value = 785.5500000000001
$ws.Range("L131").Value = 2405.5875
$ws.Range("M131").Value = 4254.45
$ws.Range("N131").Value = -12485.5875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 37626
$ws.Range("J46").Value = 37626
$ws.Range("L46").Value = 37626
$ws.Range("N46").Value = -37938
$ws.Range("H113").Value = 1658.3684
$ws.Range("I113").Value = 1700.6923
$ws.Range("J113").Value = 1566.6666
$ws.Range("K113").Value = 1700.6923
$ws.Range("L113").Value = 1566.6666
$ws.Range("M113").Value = 469.3077000000001
$ws.Range("N113").Value = -5906.6666
$ws.Range("H120").Value = 30733.334
$ws.Range("J120").Value = 30733.334
$ws.Range("L120").Value = 30733.334
$ws.Range("N120").Value = -40409.334
$ws.Range("H122").Value = 7020.8
$ws.Range("I122").Value = 2300
$ws.Range("K122").Value = 6900
$ws.Range("M122").Value = -4450
$ws.Range("H132").Value = 2834.3667
$ws.Range("I132").Value = 1383.8235
$ws.Range("J132").Value = 4731.231
$ws.Range("K132").Value = 4151.470499999999
$ws.Range("L132").Value = 14193.693
$ws.Range("M132").Value = -1621.470499999999
$ws.Range("N132").Value = -19253.693
$ws.Range("H137").Value = 38037
$ws.Range("J137").Value = 38037
$ws.Range("L137").Value = 38037
$ws.Range("N137").Value = -48237
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 39092.94
$ws.Range("J140").Value = 39092.94
$ws.Range("L140").Value = 39092.94
$ws.Range("N140").Value = -49452.94

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 4985.7144
$ws.Range("I20").Value = 3000
$ws.Range("J20").Value = 9950
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 9950
$ws.Range("M20").Value = -2774
$ws.Range("N20").Value = -10402
$ws.Range("H104").Value = 27999.6
$ws.Range("J104").Value = 27999.6
$ws.Range("L104").Value = 27999.6
$ws.Range("N104").Value = -34987.6
$ws.Range("H136").Value = 4226.56
$ws.Range("I136").Value = 1365.3334
$ws.Range("J136").Value = 6867.6924
$ws.Range("K136").Value = 4096.0002
$ws.Range("L136").Value = 20603.0772
$ws.Range("M136").Value = -1546.0002
$ws.Range("N136").Value = -25703.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 25966.666
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 37450
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 37450
$ws.Range("M33").Value = -2750
$ws.Range("N33").Value = -37950
$ws.Range("H36").Value = 25966.666
$ws.Range("I36").Value = 3000
$ws.Range("J36").Value = 37450
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 37450
$ws.Range("M36").Value = -2750
$ws.Range("N36").Value = -37950
$ws.Range("H37").Value = 41271.5
$ws.Range("I37").Value = 14999
$ws.Range("J37").Value = 50029
$ws.Range("K37").Value = 14999
$ws.Range("L37").Value = 50029
$ws.Range("M37").Value = -14796
$ws.Range("N37").Value = -50435
$ws.Range("H46").Value = 54064.453
$ws.Range("J46").Value = 54064.453
$ws.Range("L46").Value = 54064.453
$ws.Range("N46").Value = -54526.453
$ws.Range("H107").Value = 701.9
$ws.Range("I107").Value = 662.2308
$ws.Range("J107").Value = 775.5714
$ws.Range("K107").Value = 1986.6924
$ws.Range("L107").Value = 2326.7142
$ws.Range("M107").Value = -66.69240000000013
$ws.Range("N107").Value = -6166.7142
$ws.Range("H122").Value = 3820.889
$ws.Range("I122").Value = 1534.8334
$ws.Range("K122").Value = 4604.5002
$ws.Range("M122").Value = -2154.5002
$ws.Range("H126").Value = 711665.25
$ws.Range("I126").Value = 1308.6666
$ws.Range("K126").Value = 3925.9998
$ws.Range("M126").Value = -1455.9998
$ws.Range("H134").Value = 54064.453
$ws.Range("J134").Value = 54064.453
$ws.Range("L134").Value = 162193.359
$ws.Range("N134").Value = -167263.359

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N139").ClearContents()
